$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '59.947.39'
$ws.Range("E2").Value = '  +1.49%  '

$ws.Range("D3").Value = '2.307.81'
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.59%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.574'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.12%  '

$ws.Range("D9").Value = '2.305.12'
$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("E10").Value = '  +0.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.55'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.26%  '

$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.95%  '

$ws.Range("D15").Value = '59.906.06'
$ws.Range("E15").Value = '  +1.60%  '

$ws.Range("D16").Value = '2.718.16'
$ws.Range("E16").Value = '  -0.15%  '

$ws.Range("E17").Value = '  -0.76%  '

$ws.Range("D18").Value = '2.292.00'
$ws.Range("E18").Value = '  -0.73%  '

$ws.Range("E19").Value = '  -1.02%  '

$ws.Range("E20").Value = '  -2.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '312.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.42%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  -2.64%  '

$ws.Range("E28").Value = '  +4.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.89%  '

$ws.Range("B30").Value = 'SuiNetwork'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.17'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.16%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.73%  '

$ws.Range("E32").Value = '  -1.47%  '

$ws.Range("E33").Value = '  -0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.379'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.68'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.80%  '

$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.99'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '316.32'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.67%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.14%  '

$ws.Range("E42").Value = '  -0.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.29%  '

$ws.Range("E44").Value = '  -0.33%  '

$ws.Range("E45").Value = '  -1.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.568'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.19%  '

$ws.Range("E47").Value = '  +2.86%  '

$ws.Range("E48").Value = '  -0.76%  '

$ws.Range("D49").Value = '0.0₆0224'
$ws.Range("E49").Value = '  +21.47%  '

$ws.Range("E50").Value = '  +1.15%  '
